$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("NI", 2022, "SPD",    "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdniedersachsenltw22regierungsprogramm.pdf"),
    @("NI", 2022, "CDU",    "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cduniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "Grüne",  "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/gruneniedersachenltw22wahlprogrammentwurf.pdf"),
    @("NI", 2022, "FDP",    "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdpniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "AfD",    "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afdniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "FW-NI",  "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fwniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "Gesundheitsforschung", "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/gesundheitsforschungniedersachsenltw22.pdf"),
    @("NI", 2022, "Die Humanisten", "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/humanistenniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "Linke",  "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/linkeniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "DiePartei", "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/parteiniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "Piraten", "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/piratenniedersachsenltw22wahlprogramm.pdf"),
    @("NI", 2022, "Volt",   "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/voltniedersachsenltw22wahlprogramm.pdf")
)

$startRow = 303
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'TRUE"
    $ws.Cells.Item($r, 4).ClearFormats()
    $ws.Cells.Item($r, 5).Value = $row[3]
}

Write-Output "done"
